$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.212.51"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "1.905.88"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3956"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07968"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "1.879.47"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.136"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.789"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06949"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "29.237.83"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.364"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "2.117.25"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.058"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.895"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.004"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09438"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9255"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.356"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.347"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05862"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.174"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02107"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.991"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5758"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5436"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  +3.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07094"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.887"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("E49").Value = "  +6.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.069"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.02%  "
